# Update Sheets via scheduled runner
# Applies updated market-board price figures (and resulting profit
# recalculations) across several worksheets of the Leviathan_Profits workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ALC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H57").Value = 49742
$ws.Range("J57").Value = 49742
$ws.Range("L57").Value = 149226
$ws.Range("N57").Value = -150224

$ws.Range("H86").Value = 992.8333
$ws.Range("I86").Value = 992.8333
$ws.Range("K86").Value = 992.8333
$ws.Range("M86").Value = 130.1667

$ws.Range("H89").Value = 992.8333
$ws.Range("I89").Value = 992.8333
$ws.Range("K89").Value = 4964.1665
$ws.Range("M89").Value = 651.8334999999997

$ws.Range("H97").Value = 1995.7142
$ws.Range("J97").Value = 2198.8333
$ws.Range("L97").Value = 6596.499899999999
$ws.Range("N97").Value = -7588.499899999999

$ws.Range("H106").Value = 15869.125
$ws.Range("I106").Value = 2992.6
$ws.Range("K106").Value = 2992.6
$ws.Range("M106").Value = -2361.6

$ws.Range("H108").Value = 47500
$ws.Range("I108").Value = 35000
$ws.Range("J108").Value = 51666.668
$ws.Range("K108").Value = 35000
$ws.Range("L108").Value = 51666.668
$ws.Range("M108").Value = -31160
$ws.Range("N108").Value = -59346.668

$ws.Range("H110").Value = 39999
$ws.Range("J110").Value = 39999
$ws.Range("L110").Value = 39999
$ws.Range("N110").Value = -48179

$ws.Range("H135").Value = 250662.25
$ws.Range("I135").Value = 883.3333
$ws.Range("J135").Value = 999999
$ws.Range("K135").Value = 7949.9997
$ws.Range("L135").Value = 8999991
$ws.Range("M135").Value = -5414.9997
$ws.Range("N135").Value = -9005061

# ---------------------------------------------------------------------------
# ARM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 17262.828
$ws.Range("I32").Value = 4044.3699
$ws.Range("K32").Value = 4044.3699
$ws.Range("M32").Value = -3757.3699

$ws.Range("H45").Value = 4906.4053
$ws.Range("I45").Value = 5798.2915
$ws.Range("J45").Value = 3259.8462
$ws.Range("K45").Value = 5798.2915
$ws.Range("L45").Value = 3259.8462
$ws.Range("M45").Value = -5421.2915
$ws.Range("N45").Value = -4013.8462

$ws.Range("H61").Value = 2507.1428
$ws.Range("I61").Value = 2507.1428
$ws.Range("K61").Value = 2507.1428
$ws.Range("M61").Value = -2295.1428

$ws.Range("H74").Value = 1721.6285
$ws.Range("I74").Value = 1498.2693
$ws.Range("K74").Value = 1498.2693
$ws.Range("M74").Value = -624.2692999999999

$ws.Range("H77").Value = 1721.6285
$ws.Range("I77").Value = 1498.2693
$ws.Range("K77").Value = 7491.3465
$ws.Range("M77").Value = -3123.3465

$ws.Range("H136").Value = 2507.1428
$ws.Range("I136").Value = 2507.1428
$ws.Range("K136").Value = 7521.428400000001
$ws.Range("M136").Value = -4971.428400000001

# ---------------------------------------------------------------------------
# BSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H20").Value = 11142.111
$ws.Range("I20").Value = 11326.111
$ws.Range("K20").Value = 11326.111
$ws.Range("M20").Value = -11079.111

$ws.Range("H42").Value = 122999
$ws.Range("J42").Value = 122999
$ws.Range("L42").Value = 122999
$ws.Range("N42").Value = -123655

$ws.Range("H132").Value = 80845.2
$ws.Range("J132").Value = 80845.2
$ws.Range("L132").Value = 80845.2
$ws.Range("N132").Value = -90965.2

# ---------------------------------------------------------------------------
# CRP
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H107").Value = 1824.16
$ws.Range("I107").Value = 1801.9166
$ws.Range("J107").Value = 1844.6923
$ws.Range("K107").Value = 1801.9166
$ws.Range("L107").Value = 1844.6923
$ws.Range("M107").Value = 118.0834
$ws.Range("N107").Value = -5684.6923

$ws.Range("H132").Value = 2424.1143
$ws.Range("I132").Value = 2413.3076
$ws.Range("J132").Value = 2455.3333
$ws.Range("K132").Value = 7239.9228
$ws.Range("L132").Value = 7365.999899999999
$ws.Range("M132").Value = -4709.9228
$ws.Range("N132").Value = -12425.9999

# ---------------------------------------------------------------------------
# CUL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H48").Value = 146314.14
$ws.Range("I48").Value = 200839.8
$ws.Range("J48").Value = 10000
$ws.Range("K48").Value = 602519.3999999999
$ws.Range("L48").Value = 30000
$ws.Range("M48").Value = -602269.3999999999
$ws.Range("N48").Value = -30500

$ws.Range("H114").Value = 33335982
$ws.Range("I114").Value = 100000500
$ws.Range("J114").Value = 3722.75
$ws.Range("K114").Value = 300001500
$ws.Range("L114").Value = 11168.25
$ws.Range("M114").Value = -299998246
$ws.Range("N114").Value = -17676.25

$ws.Range("H122").Value = 399.0625
$ws.Range("I122").Value = 417.5
$ws.Range("J122").Value = 388
$ws.Range("K122").Value = 3757.5
$ws.Range("L122").Value = 3492
$ws.Range("M122").Value = -1307.5
$ws.Range("N122").Value = -8392

# ---------------------------------------------------------------------------
# GSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H97").Value = 17888.104
$ws.Range("I97").Value = 20894.787
$ws.Range("K97").Value = 20894.787
$ws.Range("M97").Value = -20398.787

$ws.Range("H132").Value = 8478.700000000001
$ws.Range("I132").Value = 7598.5
$ws.Range("K132").Value = 22795.5
$ws.Range("M132").Value = -20265.5

# ---------------------------------------------------------------------------
# LTW
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H7").Value = 11023.487
$ws.Range("I7").Value = 17228.666
$ws.Range("K7").Value = 17228.666
$ws.Range("M7").Value = -17116.666

$ws.Range("H22").Value = 3193.3333
$ws.Range("I22").Value = 5473.3335
$ws.Range("K22").Value = 5473.3335
$ws.Range("M22").Value = -5178.3335

$ws.Range("H27").Value = 3193.3333
$ws.Range("I27").Value = 5473.3335
$ws.Range("K27").Value = 5473.3335
$ws.Range("M27").Value = -5366.3335

$ws.Range("H40").Value = 4879
$ws.Range("I40").Value = 2422
$ws.Range("K40").Value = 2422
$ws.Range("M40").Value = -2286

$ws.Range("H61").Value = 116121.86
$ws.Range("I61").Value = 180339.28
$ws.Range("K61").Value = 180339.28
$ws.Range("M61").Value = -180137.28

$ws.Range("H74").Value = 25217
$ws.Range("J74").Value = 25217
$ws.Range("L74").Value = 25217
$ws.Range("N74").Value = -27213

$ws.Range("H77").Value = 25217
$ws.Range("J77").Value = 25217
$ws.Range("L77").Value = 75651
$ws.Range("N77").Value = -85635

$ws.Range("H82").Value = 3597.1428
$ws.Range("I82").Value = 2386.2
$ws.Range("K82").Value = 2386.2
$ws.Range("M82").Value = -2025.2

$ws.Range("H85").Value = 3597.1428
$ws.Range("I85").Value = 2386.2
$ws.Range("K85").Value = 2386.2
$ws.Range("M85").Value = -1138.2

$ws.Range("H113").Value = 116121.86
$ws.Range("I113").Value = 180339.28
$ws.Range("K113").Value = 180339.28
$ws.Range("M113").Value = -178169.28

$ws.Range("H122").Value = 5804.5557
$ws.Range("I122").Value = 6387.154
$ws.Range("K122").Value = 19161.462
$ws.Range("M122").Value = -16711.462

$ws.Range("H126").Value = 11023.487
$ws.Range("I126").Value = 17228.666
$ws.Range("K126").Value = 51685.99800000001
$ws.Range("M126").Value = -49215.99800000001

$ws.Range("H132").Value = 3956.818
$ws.Range("I132").Value = 3339.3157
$ws.Range("K132").Value = 10017.9471
$ws.Range("M132").Value = -7487.947100000001

# ---------------------------------------------------------------------------
# WVR
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H81").Value = 1321
$ws.Range("J81").Value = 1889.6
$ws.Range("L81").Value = 3779.2
$ws.Range("N81").Value = -5901.2

$ws.Range("H84").Value = 1321
$ws.Range("J84").Value = 1889.6
$ws.Range("L84").Value = 18896
$ws.Range("N84").Value = -29504

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").Value = ""

$ws.Range("H132").Value = 967415
$ws.Range("I132").Value = 1278690.1
$ws.Range("K132").Value = 3836070.3
$ws.Range("M132").Value = -3833540.3
